# Applies updated income-imputation figures (refreshed data run) to each sheet.
$wb = $excel.ActiveWorkbook

# missing_values
$ws = $wb.Worksheets.Item("missing_values")
$ws.Range("B4").Value = 186
$ws.Range("C4").Value = 1.4365152919369786
$ws.Range("B5").Value = 356
$ws.Range("C5").Value = 2.7494593759654
$ws.Range("B6").Value = 1671
$ws.Range("C6").Value = 12.905468025949954
$ws.Range("B8").Value = 2027
$ws.Range("C8").Value = 15.654927401915353
$ws.Range("B9").Value = 10658
$ws.Range("C9").Value = 82.313870868087733
$ws.Range("B10").Value = 12948
$ws.Range("B15").Value = 23404
$ws.Range("C15").Value = 374.16466826538766
$ws.Range("C16").Value = 0.11191047162270183
$ws.Range("C17").Value = 0.28776978417266186
$ws.Range("C19").Value = 0.3996802557953637
$ws.Range("B20").Value = 6160
$ws.Range("C20").Value = 98.481215027977626
$ws.Range("B21").Value = 6255
$ws.Range("C26").Value = 1.6816143497757847
$ws.Range("B27").Value = 20
$ws.Range("C27").Value = 2.2421524663677128
$ws.Range("C28").Value = 0.22421524663677131
$ws.Range("B29").Value = 22
$ws.Range("C29").Value = 2.4663677130044843
$ws.Range("B30").Value = 855
$ws.Range("C30").Value = 95.852017937219742
$ws.Range("B31").Value = 892
$ws.Range("B36").Value = 25088
$ws.Range("C36").Value = 314.26781911562068
$ws.Range("B37").Value = 57
$ws.Range("C37").Value = 0.71401728673431042
$ws.Range("B39").Value = 57
$ws.Range("C39").Value = 0.71401728673431042
$ws.Range("B40").Value = 7919
$ws.Range("C40").Value = 99.19829637980709
$ws.Range("B41").Value = 7983

# profile_missing_values
$ws = $wb.Worksheets.Item("profile_missing_values")
$ws.Range("B4").Value = 31.277750370004931
$ws.Range("B5").Value = 68.722249629995076
$ws.Range("B7").Value = 0.34533793783917122
$ws.Range("B8").Value = 12.530833744449927
$ws.Range("B9").Value = 22.792303897385299
$ws.Range("B10").Value = 23.038973852984707
$ws.Range("B11").Value = 19.980266403552047
$ws.Range("B12").Value = 14.50419338924519
$ws.Range("B13").Value = 6.8080907745436612
$ws.Range("B15").Value = 24.124321657622101
$ws.Range("B16").Value = 32.807104094721261
$ws.Range("B17").Value = 29.551060680809076
$ws.Range("B18").Value = 10.853478046373951
$ws.Range("B19").Value = 2.5160335471139614
$ws.Range("B20").Value = 0.1480019733596448
$ws.Range("B22").Value = 0.34533793783917122
$ws.Range("B23").Value = 0.1480019733596448
$ws.Range("B24").Value = 26.985693142575233
$ws.Range("B25").Value = 48.001973359644794
$ws.Range("B26").Value = 5.1307350764676869
$ws.Range("B27").Value = 14.208189442525901
$ws.Range("B28").Value = 0.39467192895905284
$ws.Range("B29").Value = 4.7853971386285155
$ws.Range("B31").Value = 4.3413912185495809
$ws.Range("B32").Value = 4.1440552540700546
$ws.Range("B33").Value = 11.001480019733597
$ws.Range("B34").Value = 73.162308830784411
$ws.Range("B35").Value = 0.39467192895905284
$ws.Range("B36").Value = 4.0453872718302915
$ws.Range("B37").Value = 2.7133695115934877
$ws.Range("B38").Value = 0.19733596447952642
$ws.Range("B40").Value = 18.697582634435125
$ws.Range("B41").Value = 0.64134188455846086
$ws.Range("B42").Value = 1.6773556980759743
$ws.Range("B43").Value = 1.1840157868771584
$ws.Range("B44").Value = 5.920078934385792
$ws.Range("B45").Value = 20.720276270350272
$ws.Range("B46").Value = 12.530833744449927
$ws.Range("B47").Value = 2.8120374938332513
$ws.Range("B48").Value = 5.8214109521460289
$ws.Range("B49").Value = 29.797730636408488
$ws.Range("B50").Value = 0.19733596447952642
$ws.Range("B52").Value = 18.894918598914654
$ws.Range("B53").Value = 81.105081401085343

# labor_incmon_imp_stochastic_reg
$ws = $wb.Worksheets.Item("labor_incmon_imp_stochastic_reg")
$ws.Range("B4").Value = 2013340.0766250594
$ws.Range("D4").Value = 454765.71875
$ws.Range("H4").Value = 2016833.1326708435
$ws.Range("J4").Value = 468985.25
$ws.Range("K4").Value = 1007619.5625
$ws.Range("L4").Value = 2606851.5
$ws.Range("M4").Value = 4042826.25

# labor_jubpenimp_stochastic_reg
$ws = $wb.Worksheets.Item("labor_jubpenimp_stochastic_reg")
$ws.Range("B4").Value = 1519927.2796204803
$ws.Range("H4").Value = 1521645.8887279662

# nonlabor_imp_stochastic_reg
$ws = $wb.Worksheets.Item("nonlabor_imp_stochastic_reg")
$ws.Range("B4").Value = 9458464.3975915201
$ws.Range("G4").Value = 927010
$ws.Range("H4").Value = 9394301.2418644987
$ws.Range("I4").Value = 120581.8125
$ws.Range("M4").Value = 930000

# labor_beneimp_stochastic_reg
$ws = $wb.Worksheets.Item("labor_beneimp_stochastic_reg")
$ws.Range("B4").Value = 868703.25247404724
$ws.Range("D4").Value = 188103.75
$ws.Range("E4").Value = 302285.875
$ws.Range("F4").Value = 860511.875
$ws.Range("G4").Value = 2000000
$ws.Range("H4").Value = 881703.5860890263
$ws.Range("J4").Value = 196485.828125
$ws.Range("K4").Value = 320000
$ws.Range("L4").Value = 860511.875
$ws.Range("M4").Value = 2000000
